# "remove excellence model from poster"
#
# The poster is a single huge slide. Removing the "excellent" model
# (predictions = slightly perturbed true values) from the write-up:
#   - deletes its row from the Results table (Table 82 / shape id 83)
#   - shrinks/re-wraps the paragraphs that mention it or sit below it
#   - shifts everything below the edited blocks upward to close the gap
#
# NOTE on numeric literals below: PowerPoint's Shape.Left/Top/Width/Height
# are single-precision (float32) underneath, so an EMU value divided by
# 12700 (EMU per point) does not always round-trip to the exact EMU we
# want -- the nearest representable float32 point value is used instead
# (at most 1 EMU away, i.e. << 1/1000 mm, from the target).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# Shape id 9  "Rounded Rectangle 8"  -- "Results" section header bar.
# Moves up (target y = 4448657 EMU).
# ---------------------------------------------------------------------
$sh9 = $s.Shapes.Item(7)
$sh9.Top = 350.2879638671875

# ---------------------------------------------------------------------
# Shape id 11 "Rounded Rectangle 10" -- "Conclusion" section header bar.
# Moves up (target x = 23595872, y = 15345569 EMU).
# ---------------------------------------------------------------------
$sh11 = $s.Shapes.Item(8)
$sh11.Left = 1857.9427490234375
$sh11.Top  = 1208.3125

# ---------------------------------------------------------------------
# Shape id 79 "TextBox 78" -- Conclusion paragraph. Reworded, grows
# taller (target x=23552944, y=16045947, cy=3416320 EMU).
# ---------------------------------------------------------------------
$sh79 = $s.Shapes.Item(57)
$sh79.Left   = 1854.5625
$sh79.Top    = 1263.46044921875
$sh79.Height = 269.0015869140625
$sh79.TextFrame.TextRange.Text = "Using techniques similar to those in prescriptive analysis (using the underlying model in addition to the point prediction) may be an effective way to gain an edge over opponents in fantasy sports, but the underlying model needs to be accurate. Our model was not accurate enough to yield a profitable system. Many improvements can be made, such as adding the predictions of a professional fantasy NBA analytics service as input to the model. Additionally, it would be useful to check how accurate a model needs to be to hit the 90% metric (e.g. what MAE leads to profitability?)"

# ---------------------------------------------------------------------
# Shape id 80 "TextBox 79" -- testing paragraph right under the
# Results header (target x=23578570, y=5119372 EMU).
# ---------------------------------------------------------------------
$sh80 = $s.Shapes.Item(58)
$sh80.Left = 1856.580322265625
$sh80.Top  = 403.10015869140625

# ---------------------------------------------------------------------
# Shape id 81 "Table 80" -- model-performance table (graphicFrame).
# Moves up (target y = 8081863 EMU).
# ---------------------------------------------------------------------
$sh81 = $s.Shapes.Item(59)
$sh81.Top = 636.3671875

# ---------------------------------------------------------------------
# Shape id 82 "TextBox 81" -- scoring-metric paragraph.
# Moves up (target y = 10043668 EMU).
# ---------------------------------------------------------------------
$sh82 = $s.Shapes.Item(60)
$sh82.Top = 790.8400268554688

# ---------------------------------------------------------------------
# Shape id 83 "Table 82" -- system-performance table (graphicFrame).
# Delete the "Randomly perturbed true values" / "X%" row (row 3, the
# "excellent" model), then re-anchor/resize the now-shorter frame
# (target x=25531672, y=12713240, cx=4726239, cy=1305248 EMU).
# ---------------------------------------------------------------------
$sh83 = $s.Shapes.Item(61)
$tbl83 = $sh83.Table
$tbl83.Rows.Item(3).Delete()

$sh83.Left   = 2010.367919921875
$sh83.Top    = 1001.0425415039062
$sh83.Width  = 372.1448059082031
$sh83.Height = 102.77543640136719

# ---------------------------------------------------------------------
# Shape id 95 "TextBox 94" -- Optimization intro paragraph. Drops the
# sentence about the "excellent" model and shrinks (target cy=1938992
# EMU).
# ---------------------------------------------------------------------
$sh95 = $s.Shapes.Item(68)
$sh95.Height = 152.67654418945312
$sh95.TextFrame.TextRange.Text = "In addition to the base optimization model (maximize score and variance), we created and tested other systems. One was a baseline version using point-predictions only. Another used the variances to adjust the point-predictions, without explicitly using the variance in the optimization model."

# ---------------------------------------------------------------------
# Shape id 96 "TextBox 95" -- "Table 1" caption.
# Moves up (target y = 7763330 EMU).
# ---------------------------------------------------------------------
$sh96 = $s.Shapes.Item(69)
$sh96.Top = 611.2858276367188

# ---------------------------------------------------------------------
# Shape id 97 "Rectangle 96" -- footnote about ORTs.
# Moves up (target y = 9741998 EMU).
# ---------------------------------------------------------------------
$sh97 = $s.Shapes.Item(70)
$sh97.Top = 767.0864868164062

# ---------------------------------------------------------------------
# Shape id 98 "TextBox 97" -- closing paragraph about the 90%
# profitability mark. Reworded, shrinks (target y=14081864,
# cy=1200329 EMU).
# ---------------------------------------------------------------------
$sh98 = $s.Shapes.Item(71)
$sh98.Top    = 1108.8082275390625
$sh98.Height = 94.51409912109375
$sh98.TextFrame.TextRange.Text = "One can see that we were unable to reach the 90% profitability mark in any of our three designs. This may be due to our initial predictions being off by ~25%, on average."
